# Apply updated dSF (column F) values to Sheet1
# repull data, push all data, mean calculation

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    6  = 1
    7  = 0
    10 = 0
    19 = 4
    21 = 4
    22 = -3
    27 = -3
    28 = 0
    32 = -3
    33 = -5
    36 = -1
    38 = -5
    40 = 4
    41 = 1
    43 = -1
    47 = -4
    50 = 4
    53 = -2
    59 = -1
    62 = 3
    67 = -3
    68 = 2
    70 = -6
    73 = -5
    74 = -2
    75 = -2
    80 = -7
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
